# Included sql property file and excel sheet for sql row and column count.
#
# Adds two new columns to the "addVisitor" sheet - sqlRecordCount and
# sqlColCount - together with a data row holding the SQL record count
# (252) and column count (5), and removes the old leftover sample rows
# (the stale e-mail addresses that used to live in A6/A7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the old leftover rows (3-7), keeping only the header row and
#     a fresh row 2 for the new data ---
$ws.Range("A3:A7").EntireRow.Delete()

# --- Header row (row 1): add the two new headers ---
$ws.Range("B1").Value = "sqlRecordCount"
$ws.Range("C1").Value = "sqlColCount"

# --- Data row (row 2) ---
$ws.Range("A2").Value = 10

# sqlRecordCount value - stored as text "252"
$b2 = $ws.Range("B2")
$b2.NumberFormat = "@"
$b2.Value = "252"
$b2.NumberFormat = "General"

# sqlColCount value - stored as text "5"
$c2 = $ws.Range("C2")
$c2.NumberFormat = "@"
$c2.Value = "5"
$c2.NumberFormat = "General"

# --- Widen the new sqlRecordCount column so the header fits ---
$ws.Columns.Item(2).ColumnWidth = 13.5
